# Loading of xlsx file: add a new "Format" column (F) to the PPE structure
# sheet, and populate the "ddmmrrrr" format value for the two date fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "Format", formatted like the rest of the header row (E1)
$ws.Range("F1").Value = "Format"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Row 5 ("dátum vytvorenia súboru – ddmmrrrr" / fileCreated) gets its format
$ws.Range("F5").Value = "ddmmrrrr"

# Row 10 ("dátum služby Vyplaťte dňa VD–ddmmrrrr" / payOutDate) gets its format
$ws.Range("F10").Value = "ddmmrrrr"

# Move/save the selection to A11, matching where the cursor ended up after edits
$ws.Range("A11").Select()
